# issue #5: add legislator_id, name, date into dataframe
# Add three new columns (date, legislator_name, legislator_id) to the
# "股票" (stocks) worksheet, filling every existing data row with the
# filing date, legislator name, and legislator id.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Header row - copy the existing header formatting (bold + border) onto the
# new header cells, then fill in the labels.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# Determine how many data rows already exist (rows below the header)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $dateCell = $ws.Cells.Item($r, 8)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2012-04-02"
    $ws.Cells.Item($r, 9).Value = "丁守中"
    $ws.Cells.Item($r, 10).Value = 515

    # Copy the existing data-row formatting (border etc.) onto the new cells,
    # applied after the values so the paste doesn't get clobbered by the
    # value assignment.
    $ws.Range("G$r").Copy() | Out-Null
    $ws.Range("H$r:J$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
